$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: the "Nom" (C2) / "Prenom" (D2) cells of the grade sheet still held
# the generic placeholder values "Name" / "Forename" instead of the actual
# sender's identity that is mail-merged into the report. Replace them with
# the real values so the generated report no longer leaks the placeholder.
$ws.Range("C2").Value = "REY"
$ws.Range("D2").Value = "Alexandre"

# Leave the selection where the author ended up after making the edit.
$ws.Range("E16").Select()
